$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46018
$ws.Range("B2").Value = 92.34
$ws.Range("C2").Value = 88.08
$ws.Range("D2").Value = 80.53
$ws.Range("E2").Value = 70.95999999999999
$ws.Range("F2").Value = 69.36
$ws.Range("G2").Value = 68.95
$ws.Range("H2").Value = 69.17
$ws.Range("I2").Value = 72.51000000000001
$ws.Range("J2").Value = 86.79000000000001
$ws.Range("K2").Value = 88.12
$ws.Range("L2").Value = 75.37
$ws.Range("M2").Value = 69.26000000000001
$ws.Range("N2").Value = 68.84999999999999
$ws.Range("O2").Value = 69.14
$ws.Range("P2").Value = 69.36
$ws.Range("Q2").Value = 69.98999999999999
$ws.Range("R2").Value = 81.98
$ws.Range("S2").Value = 91.16
$ws.Range("T2").Value = 96.56999999999999
$ws.Range("U2").Value = 98.40000000000001
$ws.Range("V2").Value = 91.77
$ws.Range("W2").Value = 86.70999999999999
$ws.Range("X2").Value = 86.63
$ws.Range("Y2").Value = 80.81999999999999
$ws.Range("Z2").Value = 80.12
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 92.03
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 97.48
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 90.20999999999999
$ws.Range("AG2").Value = "3h-15h"
